# Rename the "_old"/"_new" column header suffixes to the respective
# input-file format-version suffixes ("_FV2404" / "_FV2410"), turn the
# header+data range into a real Excel Table ("Table1"), and freeze the
# header row - matching the upstream commit "Use <formatversion> as
# suffix for table headers" / "Adjust xlsx export to new header
# formatting".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $ws.Range("A1:U1")

# Stash the header's existing formatting (bold font, grey fill, border,
# centered+wrapped) in a scratch row so it survives table creation.
$stash = $ws.Range("A200:U200")
$headerRange.Copy()
$stash.PasteSpecial(-4122)   # xlPasteFormats

# Reset the header to the default style before adding the table: this
# engine auto-captures any pre-existing header formatting into a new
# headerRowDxfId, which the target workbook does not have.
$headerRange.Style = "Normal"

$rng = $ws.Range("A1:U75")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Put the original header formatting back and clean up the scratch row.
$stash.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$stash.Clear()
$excel.CutCopyMode = $false

# Freeze the header row (pane split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
